$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$fcs = $ws.Range("G48:G52").FormatConditions
$n = $fcs.Count()
for ($i=$n; $i -ge 1; $i--) {
  $fc = $fcs.Item($i)
  $fc.Delete()
}
